$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the refreshed coinranking.com snapshot values.
# A few of the new "Price" values are plain decimal-looking strings (e.g. "7.80",
# "1.00"). Assigning those straight to .Value would make Excel auto-convert them to
# numbers (dropping the literal trailing zeros), so for those cells we temporarily
# force a Text number format, assign the literal string, then restore the original
# (Normal) cell style so no visible formatting is left behind.
$ws.Range("D2").Value = '60.299.05'
$ws.Range("E2").Value = '  -3.75%  '
$ws.Range("D3").Value = '3.306.43'
$ws.Range("E3").Value = '  -4.13%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '560.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.69'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.31%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '3.308.60'
$ws.Range("E8").Value = '  -4.10%  '
$ws.Range("E9").Value = '  -0.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.80'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.16%  '
$ws.Range("E11").Value = '  -3.30%  '
$ws.Range("E12").Value = '  -1.40%  '
$ws.Range("D13").Value = '3.879.29'
$ws.Range("E13").Value = '  -3.98%  '
$ws.Range("E14").Value = '  +0.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.22'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.93%  '
$ws.Range("D16").Value = '3.309.89'
$ws.Range("E16").Value = '  -3.85%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000166'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.40%  '
$ws.Range("D18").Value = '60.342.88'
$ws.Range("E18").Value = '  -3.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.16'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.39'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.62'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '373.83'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.50%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '74.16'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.550'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.79%  '
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("D26").Value = '3.471.39'
$ws.Range("E26").Value = '  -3.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000105'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.30%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.174'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.997'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.22'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.95%  '
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.05'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.03%  '
$ws.Range("E33").Value = '  -4.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '22.59'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.28'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.17'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.00%  '
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.76'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.70%  '
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '165.94'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.96%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.53'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '27.88'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -13.18%  '
$ws.Range("D41").Value = '3.343.58'
$ws.Range("E41").Value = '  -4.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0737'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.93'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.752'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.22'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.75%  '
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.60'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.00%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.12'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.10%  '
$ws.Range("D48").Value = '2.380.66'
$ws.Range("E48").Value = '  -7.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.59'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.69'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.74%  '
